$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Re-apply the default style to the existing data range so it gets a
# dedicated (but visually identical) cell style entry.
$sheet1.Range("A1:E5").Style = "Normal"

# Add the new worksheet after Sheet1 and populate it with data.
$sheet2 = $wb.Worksheets.Add($null, $sheet1)
$sheet2.Name = "Sheet2"

$data = @(
    @("new a1", "new b1", "new c1", "new d1", "new e1"),
    @("new a2", "new b2", "new c2", "new d2", "new e2"),
    @("new a3", "new b3", "new c3", "new d3", "new e3"),
    @("new a4", "new b4", "new c4", "new d4", "new e4"),
    @("new a5", "new b5", "new c5", "new d5", "new e5")
)

for ($r = 0; $r -lt 5; $r++) {
    for ($c = 0; $c -lt 5; $c++) {
        $sheet2.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

$sheet2.Range("A1:E5").Style = "Normal"

$sheet2.Activate()
